# Update "想去人数" (number of people interested) counts across the
# 展览, 演出, 本地生活 and 全部类型 worksheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 25
$ws1.Range("F9").Value = 1080
$ws1.Range("F13").Value = 1361
$ws1.Range("F21").Value = 645
$ws1.Range("F26").Value = 5220
$ws1.Range("F29").Value = 2434
$ws1.Range("F30").Value = 5839
$ws1.Range("F31").Value = 125
$ws1.Range("F40").Value = 679

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 7
$ws2.Range("F34").Value = 4

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value = 382
$ws3.Range("F7").Value = 220

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 382
$ws4.Range("F7").Value = 220
$ws4.Range("F8").Value = 220
$ws4.Range("F13").Value = 25
$ws4.Range("F16").Value = 1080
$ws4.Range("F21").Value = 1361
$ws4.Range("F29").Value = 645
$ws4.Range("F34").Value = 5220
$ws4.Range("F37").Value = 2434
$ws4.Range("F38").Value = 5839
$ws4.Range("F39").Value = 125
$ws4.Range("F46").Value = 679
